# Apply the edit: take the absolute value of the "volume" column (K),
# and recompute the "compactness" column (L) using the corrected formula
# L = area^3 / (36 * pi * volume^2), using the corrected (absolute) volume.
#
# Columns: J = area (10), K = volume (11), L = compactness (12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pi = 3.14159265358979323846

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $area = $ws.Cells.Item($r, 10).Value2
    $volume = $ws.Cells.Item($r, 11).Value2

    $volumeAbs = [Math]::Abs($volume)
    $compactness = ($area * $area * $area) / (36 * $pi * $volumeAbs * $volumeAbs)

    $ws.Cells.Item($r, 11).Value = $volumeAbs
    $ws.Cells.Item($r, 12).Value = $compactness
}
